# Update "想去人数" (interest count) figures in column F across the
# "展览" (Exhibitions) and "全部类型" (All Types) sheets, matching
# the refreshed output data.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 15040
$wsExhibit.Range("F3").Value = 19142
$wsExhibit.Range("F4").Value = 147
$wsExhibit.Range("F5").Value = 142
$wsExhibit.Range("F14").Value = 166
$wsExhibit.Range("F15").Value = 228
$wsExhibit.Range("F16").Value = 67
$wsExhibit.Range("F17").Value = 1477
$wsExhibit.Range("F20").Value = 98
$wsExhibit.Range("F21").Value = 238
$wsExhibit.Range("F22").Value = 7993
$wsExhibit.Range("F26").Value = 66
$wsExhibit.Range("F29").Value = 6071
$wsExhibit.Range("F32").Value = 171
$wsExhibit.Range("F34").Value = 290
$wsExhibit.Range("F35").Value = 5461
$wsExhibit.Range("F36").Value = 671
$wsExhibit.Range("F37").Value = 16
$wsExhibit.Range("F39").Value = 50

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 15040
$wsAll.Range("F3").Value = 19142
$wsAll.Range("F4").Value = 147
$wsAll.Range("F5").Value = 142
$wsAll.Range("F14").Value = 166
$wsAll.Range("F15").Value = 228
$wsAll.Range("F16").Value = 67
$wsAll.Range("F17").Value = 1477
$wsAll.Range("F21").Value = 98
$wsAll.Range("F22").Value = 238
$wsAll.Range("F23").Value = 7993
$wsAll.Range("F27").Value = 66
$wsAll.Range("F32").Value = 6071
$wsAll.Range("F35").Value = 171
$wsAll.Range("F37").Value = 290
$wsAll.Range("F38").Value = 5461
$wsAll.Range("F39").Value = 671
$wsAll.Range("F40").Value = 16
$wsAll.Range("F42").Value = 50
